$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 84

# Column A holds a date written as plain text (matches the rest of the
# column, e.g. "2025/10/09"); a leading apostrophe forces Excel to store
# it as literal text instead of auto-converting it to a date serial.
$ws.Cells.Item($newRow, 1).Value = "'2025/10/09"
$ws.Cells.Item($newRow, 2).Value = "木"
$ws.Cells.Item($newRow, 3).Value = 17
$ws.Cells.Item($newRow, 4).Value = 176
